# The workbook gets a brand-new weekly price record inserted as row 22
# (pushing the existing rows 22-142 down to 23-143).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22; everything below shifts down by one.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new record's data.
$ws.Range("A22").Value2 = 6
$ws.Range("B22").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C22").Value2 = "Metropolitana"
$ws.Range("D22").Value2 = 45145
$ws.Range("E22").Value2 = 13
$ws.Range("F22").Value2 = 100114007
$ws.Range("G22").Value2 = "Jengibre"
$ws.Range("H22").Value2 = "Sin especificar"
$ws.Range("I22").Value2 = "Primera"
$ws.Range("J22").Value2 = 300
$ws.Range("K22").Value2 = 16000
$ws.Range("L22").Value2 = 17000
$ws.Range("M22").Value2 = 16500
$ws.Range("N22").Value2 = "`$/caja 13 kilos"
$ws.Range("O22").Value2 = "Perú"
$ws.Range("P22").Value2 = 1269
$ws.Range("Q22").Value2 = 13
$ws.Range("R22").Value2 = "Hortaliza"
